# Apply cell value updates to match the target cryptos.xlsx revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.462.41"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "3.403.14"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'596.34"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'141.49"
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.403.35"
$ws.Range("E7").Value = "  -2.36%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D10").Value = "'7.92"
$ws.Range("E10").Value = "  +4.69%  "
$ws.Range("E11").Value = "  -6.24%  "
$ws.Range("D12").Value = "'0.404"
$ws.Range("E12").Value = "  -4.86%  "
$ws.Range("D13").Value = "3.978.93"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("E14").Value = "  -6.92%  "
$ws.Range("D15").Value = "'29.35"
$ws.Range("E15").Value = "  -6.93%  "
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "65.531.79"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.393.97"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "'10.30"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "'6.10"
$ws.Range("E20").Value = "  -5.78%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("D22").Value = "'412.54"
$ws.Range("E22").Value = "  -6.07%  "
$ws.Range("E23").Value = "  -5.80%  "
$ws.Range("D24").Value = "'76.88"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "3.539.77"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("E27").Value = "  -9.55%  "
$ws.Range("E28").Value = "  -6.39%  "
$ws.Range("E29").Value = "  -7.68%  "
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  -5.76%  "
$ws.Range("D33").Value = "'1.44"
$ws.Range("E33").Value = "  -9.72%  "
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "3.398.82"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.67"
$ws.Range("E37").Value = "  -7.38%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.50"
$ws.Range("E38").Value = "  -9.19%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'7.48"
$ws.Range("E40").Value = "  -5.73%  "
$ws.Range("D41").Value = "'167.21"
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("E42").Value = "  -4.66%  "
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("E44").Value = "  -7.69%  "
$ws.Range("E45").Value = "  -11.04%  "
$ws.Range("D46").Value = "'45.33"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "'26.32"
$ws.Range("E47").Value = "  -9.62%  "
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "'7.02"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  -7.17%  "
$ws.Range("E51").Value = "  -7.15%  "
